{"js": "// Apply the \"BusinessTripCostComponent\" -> \"WorkArriveDepartPermit\" rename,\n// the description text rewrite, the revision-date fix, the refreshed JWT\n// sample token, and the recordID sample value update (with the \"_GoBack\"\n// bookmark following the edit point), matching the authored diff.\n\nconst body = context.document.body;\n\n// 1) Rename every \"BusinessTripCostComponent\" occurrence (title, method\n//    name, file path/name, JSON keys, sample code) to \"WorkArriveDepartPermit\".\nconst apiNameResults = body.search(\"BusinessTripCostComponent\", { matchCase: true });\napiNameResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < apiNameResults.items.length; i++) {\n  apiNameResults.items[i].insertText(\"WorkArriveDepartPermit\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Replace the Indonesian description with the new text.\nconst descResults = body.search(\"Menghapus Data Komponen Biaya Perjalanan Bisnis\", { matchCase: true });\ndescResults.load(\"items\");\nawait context.sync();\nif (descResults.items.length > 0) {\n  descResults.items[0].insertText(\n    \"Menghapus Data Izin Pulang Cepat dan Datang Terlambat\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 3) Fix the revision date from the 23rd to the 25th of November.\nconst dateResults = body.search(\"23th\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  const digitResults = dateResults.items[0].search(\"3\", { matchCase: true });\n  digitResults.load(\"items\");\n  await context.sync();\n  if (digitResults.items.length > 0) {\n    digitResults.items[0].insertText(\"5\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 4) Refresh the sample JWT token in the sample request code block.\nconst oldJwt =\n  \"eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJsb2dnZWRJbkFzIjoidGVndWgucHJhdGFtYSIsImlhdCI6MTYwNjA5Nzg4MH0.d1AB_XF31WOFS7dhxvEHyJmPybR5ju4YHiuF_ZbSf5Q\";\nconst newJwt =\n  \"eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJsb2dnZWRJbkFzIjoidGVndWgucHJhdGFtYSIsImlhdCI6MTYwNjI2OTA1NH0.NjJJegg6WRVQ3LHksbKcni92MkyzjfYpxzrFvgLu2FQ\";\nconst jwtResults = body.search(oldJwt, { matchCase: true });\njwtResults.load(\"items\");\nawait context.sync();\nif (jwtResults.items.length > 0) {\n  jwtResults.items[0].insertText(newJwt, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 5) The stray \"_GoBack\" bookmark (Word's \"last edit\" marker) moves to\n//    track the new edit point: drop the old one and drop a fresh one right\n//    before the updated sample recordID value.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst recordIdResults = body.search(\"810000000000001\", { matchCase: true });\nrecordIdResults.load(\"items\");\nawait context.sync();\nif (recordIdResults.items.length > 0) {\n  const oldValueResults = recordIdResults.items[0].search(\"81\", { matchCase: true });\n  oldValueResults.load(\"items\");\n  await context.sync();\n  if (oldValueResults.items.length > 0) {\n    const target = oldValueResults.items[0];\n    const startRange = target.getRange(\"Start\");\n    startRange.insertBookmark(\"_GoBack\");\n    target.insertText(\"40\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Apply the \"BusinessTripCostComponent\" -> \"WorkArriveDepartPermit\" rename,\n# the description text rewrite, the revision-date fix, the refreshed JWT\n# sample token, and the recordID sample value update (with the \"_GoBack\"\n# bookmark following the edit point), matching the authored diff.\n\n$d = $word.ActiveDocument\n\n# 1) Rename every \"BusinessTripCostComponent\" occurrence (title, method\n#    name, file path/name, JSON keys, sample code) to \"WorkArriveDepartPermit\".\n$find = $d.Content.Find\n$find.Execute(\"BusinessTripCostComponent\", $false, $false, $false, $false, $false, $true, 1, $false, \"WorkArriveDepartPermit\", 2)\n\n# 2) Replace the Indonesian description with the new text.\n$find2 = $d.Content.Find\n$find2.Execute(\"Menghapus Data Komponen Biaya Perjalanan Bisnis\", $false, $false, $false, $false, $false, $true, 1, $false, \"Menghapus Data Izin Pulang Cepat dan Datang Terlambat\", 2)\n\n# 3) Fix the revision date from the 23rd to the 25th of November (only the\n#    leading digits are touched so the superscript \"th\" keeps its formatting).\n$find3 = $d.Content.Find\n$find3.Execute(\"23\", $false, $false, $false, $false, $false, $true, 1, $false, \"25\", 2)\n\n# 4) Refresh the sample JWT token in the sample request code block.\n$oldJwt = \"eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJsb2dnZWRJbkFzIjoidGVndWgucHJhdGFtYSIsImlhdCI6MTYwNjA5Nzg4MH0.d1AB_XF31WOFS7dhxvEHyJmPybR5ju4YHiuF_ZbSf5Q\"\n$newJwt = \"eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJsb2dnZWRJbkFzIjoidGVndWgucHJhdGFtYSIsImlhdCI6MTYwNjI2OTA1NH0.NjJJegg6WRVQ3LHksbKcni92MkyzjfYpxzrFvgLu2FQ\"\n$find4 = $d.Content.Find\n$find4.Execute($oldJwt, $false, $false, $false, $false, $false, $true, 1, $false, $newJwt, 2)\n\n# 5) The stray \"_GoBack\" bookmark (Word's \"last edit\" marker) moves to track\n#    the new edit point: re-adding it collapsed right before the updated\n#    sample recordID value removes it from its old spot and drops it here.\n$bookmarkRange = $d.Content\n$bookmarkRange.Find.Execute(\"810000000000001\", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)\n$bookmarkRange.End = $bookmarkRange.Start\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n# 6) Update the sample recordID value from 81... to 40...\n$valueRange = $d.Content\n$valueRange.Find.Execute(\"810000000000001\", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)\n$valueRange.End = $valueRange.Start + 2\n$valueRange.Text = \"40\"\n"}
